# Apply the LOB1254.xlsx content update:
#  - Row 10 (Objetivos): B/C text replaced with the docente name.
#  - Row 13 gains an "A" label ("Programa resumido:") and B/C become "Semestral".
#  - Rows 14-16 shift to the next section's content (Short syllabus / Programa / Syllabus).
#  - Row 17 becomes just "Avaliação:" (B/C cleared).
#  - Row 18 gains B/C with the docente name again; rows 19-21 shift up one slot.
#  - Row 22 (old Bibliografia content row) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos ---
$ws.Range("B10").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C10").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Rows.Item(10).RowHeight = 60

# --- Row 13: gains column A + new B/C ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: Short syllabus ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Endogenous and exogenous processes of the Earth. Materials constituting the earth's crust (minerals and rocks)."
$ws.Range("C14").Value = "Endogenous and exogenous processes of the Earth. Materials constituting the earth's crust (minerals and rocks)."
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: Programa: / 01/01/2022 ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: Syllabus ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Brief history of geology. Materials constituting the earth's crust (minerals and rocks). Origin and constitution of the universe, the solar system and the earth. Internal structure of the earth. Composition of the earth. Endogenous and exogenous processes (internal and external dynamics of the earth). Plate tectonics theory.  Igneous rocks and vulcanismo. Metamorphic rocks and metamorphism. Sedimentary rocks. Weathering, erosion, sediment transport. Geological structure. Geological time and stratigraphy."
$ws.Range("C16").Value = "Brief history of geology. Materials constituting the earth's crust (minerals and rocks). Origin and constitution of the universe, the solar system and the earth. Internal structure of the earth. Composition of the earth. Endogenous and exogenous processes (internal and external dynamics of the earth). Plate tectonics theory.  Igneous rocks and vulcanismo. Metamorphic rocks and metamorphism. Sedimentary rocks. Weathering, erosion, sediment transport. Geological structure. Geological time and stratigraphy."
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: Avaliação: only (clear B/C) ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).RowHeight = 15

# --- Row 18: Método: + docente name ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C18").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: Critério: + método text ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: Norma de recuperação: + critério text ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de provas  e atividades."
$ws.Range("C20").Value = "Média ponderada de provas  e atividades."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: Bibliografia: + norma de recuperação text ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "1 (uma) prova escrita"
$ws.Range("C21").Value = "1 (uma) prova escrita"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: removed entirely (old Bibliografia long-text row) ---
$ws.Rows.Item(22).Delete()
